$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E to be treated as plain text so Excel does not
# reinterpret dotted/percent strings as numbers, dates, or scientific notation.
$ws.Range("B2:E51").NumberFormat = "@"

$data = @(
    @('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '68.942.04', '  +0.69%  '),
    @('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.475.33', '  +0.67%  '),
    @('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.05%  '),
    @('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '561.85', '  +0.40%  '),
    @('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '164.52', '  +0.41%  '),
    @('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.12%  '),
    @('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.511', '  +1.56%  '),
    @('LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '2.472.16', '  +0.55%  '),
    @('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.158', '  +5.36%  '),
    @('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.165', '  +0.66%  '),
    @('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.334', '  -0.90%  '),
    @('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '4.85', '  +0.83%  '),
    @('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '68.867.85', '  +0.77%  '),
    @('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000172', '  +1.59%  '),
    @('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '23.70', '  +1.88%  '),
    @('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '10.68', '  -2.29%  '),
    @('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '340.03', '  -0.52%  '),
    @('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.95', '  -3.15%  '),
    @('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '3.81', '  +1.03%  '),
    @('SuiNetwork', 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui', '1.89', '  +1.24%  '),
    @('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.04%  '),
    @('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '67.07', '  -0.57%  '),
    @('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '3.68', '  -0.44%  '),
    @('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.29', '  +2.24%  '),
    @('PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0831', '  +0.12%  '),
    @('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '7.26', '  +0.82%  '),
    @('FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  +0.06%  '),
    @('Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '430.89', '  -0.06%  '),
    @('Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.15', '  -1.42%  '),
    @('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.64', '  -1.15%  '),
    @('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '159.24', '  +1.56%  '),
    @('WhiteBITCoin', 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt', '19.01', '  +0.00%  '),
    @('USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '1.00', '  -0.02%  '),
    @('Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.108', '  -2.16%  '),
    @('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.90', '  +0.25%  '),
    @('RenderToken', 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render', '4.45', '  +0.09%  '),
    @('PolygonEcosystemToken', 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol', '0.299', '  -2.10%  '),
    @('Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.49', '  -1.90%  '),
    @('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.08', '  -0.14%  '),
    @('dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '2.07', '  -0.07%  '),
    @('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.40', '  +1.52%  '),
    @('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '131.23', '  -2.44%  '),
    @('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.0722', '  +0.88%  '),
    @('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '0.488', '  +1.55%  '),
    @('Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.567', '  +0.56%  '),
    @('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.0918', '  +1.27%  '),
    @('BitgetToken', 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb', '1.12', '  +0.24%  '),
    @('Optimism', 'https://coinranking.com/coin/n1p-s_gm1+optimism-op', '1.39', '  -2.17%  '),
    @('THORChain', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', '5.00', '  -5.11%  '),
    @('InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '17.00', '  -2.68%  '),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}